# Fruta / hortaliza, semanal
#
# A new weekly price record is inserted as row 25 (Castle Brite / Primera,
# Region Metropolitana, 2022-11-30), which pushes every existing record
# from the old row 25 through row 45 down by one row (old row N becomes
# new row N+1). The sheet's used range grows from A1:T45 to A1:T46.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 25; Excel shifts rows 25:45 down to 26:46
# and grows the sheet's dimension automatically.
$ws.Rows.Item(25).Insert()

# Populate the newly inserted row 25 with the new weekly record.
$ws.Cells.Item(25, 1).Value = 4
$ws.Cells.Item(25, 2).Value = 'Feria Lagunitas de Puerto Montt'
$ws.Cells.Item(25, 3).Value = 'Los Lagos'
$ws.Cells.Item(25, 4).Value = '2022-11-30'
$ws.Cells.Item(25, 5).Value = 10
$ws.Cells.Item(25, 6).Value = 'Fruta'
$ws.Cells.Item(25, 7).Value = 100103
$ws.Cells.Item(25, 8).Value = 'Frutos de hueso (carozo)'
$ws.Cells.Item(25, 9).Value = 100103003
$ws.Cells.Item(25, 10).Value = 'Damasco'
$ws.Cells.Item(25, 11).Value = 'Castle Brite'
$ws.Cells.Item(25, 12).Value = 'Primera'
$ws.Cells.Item(25, 13).Value = 400
$ws.Cells.Item(25, 14).Value = 23000
$ws.Cells.Item(25, 15).Value = 24000
$ws.Cells.Item(25, 16).Value = 23500
$ws.Cells.Item(25, 17).Value = '$/caja 16 kilos'
$ws.Cells.Item(25, 18).Value = 'Región Metropolitana'
$ws.Cells.Item(25, 19).Value = 1469
$ws.Cells.Item(25, 20).Value = 16
